$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 370
$ws1.Range("F4").Value = 285
$ws1.Range("F5").Value = 4181

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 7

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 370
$ws4.Range("F4").Value = 285
$ws4.Range("F5").Value = 4181
$ws4.Range("F7").Value = 7
